# 868.docx edit: convert a Heading1 title + bold byline paragraph into
# pandoc-style Title/Authors paragraphs with one run per word/space token,
# and drop the bookmark that used to wrap the old heading.

$d = $word.ActiveDocument

# --- 1. Remove the old paragraphs (Heading1 title, bold "By Dorothy Day") ---
$d.Paragraphs.Item(2).Range.Delete()   # "By Dorothy Day"
$d.Paragraphs.Item(1).Range.Delete()   # old Heading1 title

# --- 2. Strip the now-empty bookmarkStart/bookmarkEnd pair that wrapped the
#        old title paragraph. Deleting a zero-length range at the very start
#        of the document clears one dangling bookmark marker at a time
#        without touching real content, so do it twice (start, then end). ---
$d.Range(0, 0).Delete()
$d.Range(0, 0).Delete()

# --- 3. Insert the replacement paragraphs as literal OOXML so each word and
#        space becomes its own <w:r> (matching the pandoc-generated target),
#        instead of collapsing into a single merged run. ---
$titleAndAuthorsXml = '<w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r><w:t xml:space="preserve">Our</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Country</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Passes</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">from</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Undeclared</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">War</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">to</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Declared</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">War</w:t></w:r><w:r><w:t xml:space="preserve">;</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">We</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Continue</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Our</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Christian</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Pacifist</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Stand</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Authors"/></w:pPr><w:r><w:t xml:space="preserve">Dorothy</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Day</w:t></w:r></w:p>'

$null = $d.Range(0, 0).InsertXML($titleAndAuthorsXml)
